$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-20 Thursday" "2025-11-21 Friday"

Replace-Text "942×3=" "374×9="
Replace-Text "911×5=" "949×8="
Replace-Text "176×2=" "533×5="
Replace-Text "784×5=" "749×4="
Replace-Text "791×2=" "740×7="

Replace-Text "850×5=" "569×5="
Replace-Text "518×5=" "940×6="
Replace-Text "112×5=" "475×4="
Replace-Text "633×7=" "108×2="
Replace-Text "630×2=" "782×9="

Replace-Text "406×3=" "488×8="
Replace-Text "184×9=" "925×8="
Replace-Text "815×3=" "550×7="
Replace-Text "306×7=" "886×7="
Replace-Text "418×7=" "829×3="

Replace-Text "527×2=" "600×2="
Replace-Text "377×4=" "450×2="
Replace-Text "268×5=" "471×5="
Replace-Text "308×6=" "202×3="
Replace-Text "186×2=" "620×3="

Replace-Text "337×2=" "980×7="
Replace-Text "598×9=" "266×8="
Replace-Text "429×7=" "934×9="
Replace-Text "407×2=" "958×4="
Replace-Text "339×9=" "667×5="
